$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DID numbers (duplication prevention)
$ws.Range("A2").Value = 448455642999
$ws.Range("A7").Value = 448455642998

# Update selection to D15 (upload to service leaves selection elsewhere)
$ws.Range("D15").Select()
